$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.567.36"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.337.75"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.23"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.31"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.580"
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.336.52"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.95"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.759.39"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.495.11"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.335.02"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.69"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "315.03"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.61"
$ws.Range("E22").Value = "  -2.67%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.97"
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("E28").Value = "  +3.06%  "
$ws.Range("E29").Value = "  +8.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.74"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.13"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0735"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.385"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.06"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.13"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "329.48"
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.08"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "138.05"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.53"
$ws.Range("E44").Value = "  +1.85%  "
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.29"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.568"
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0223"
$ws.Range("E49").Value = "  +6.50%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0216"
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.95"
$ws.Range("E51").Value = "  -0.76%  "
